# update CarList to use RTK Query
# Adds 8 new invoice rows (IDs 49-56) to Sheet1 ("CarList"), one of which
# (row 20) introduces a new client name "AK AK" alongside the existing
# "Adam Testowy" client used by the other new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function SetRow($Row, $Id, $Client, $Paid, $Due, $CreatedDate, $PaymentDate) {
    $ws.Range("A$Row").Value = $Id
    $ws.Range("B$Row").Value = $Client
    $ws.Range("C$Row").Value = $Paid
    $ws.Range("D$Row").Value = $Due

    $ws.Range("E$Row").Value = $CreatedDate
    $ws.Range("E$Row").NumberFormat = "dd.MM.yyyy"

    $ws.Range("F$Row").Value = $PaymentDate
    $ws.Range("F$Row").NumberFormat = "dd.MM.yyyy"
}

SetRow 15 49 "Adam Testowy" 420 420 45331 45331
SetRow 16 50 "Adam Testowy" 750 750 45333.488132211896 45333.488132212464
SetRow 17 51 "Adam Testowy" 500 500 45333.50467125974 45333.50467125977
SetRow 18 52 "Adam Testowy" 250 250 45333.52913617004 45333.529136170066
SetRow 19 53 "Adam Testowy" 450 450 45337.634505044945 45337.634505045404
SetRow 20 54 "AK AK" 720 720 45337.75045813944 45337.7504581405
SetRow 21 55 "Adam Testowy" 960 960 45337.754077653444 45337.75407765347
SetRow 22 56 "Adam Testowy" 420 420 45413.52797942728 45413.527979429186
